$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# Every cell that currently holds the "Ready for handoff" shared string needs
# to be rewritten so no cell references the old text; the engine drops the
# now-unused shared-string entry and reuses one for all the new occurrences.

$ws1 = $wb.Worksheets.Item(1)   # "Overview"
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

$ws2 = $wb.Worksheets.Item(2)   # "zh-cn"
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

$ws3 = $wb.Worksheets.Item(3)   # "de-de"
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# --- Column width changes ---
# Target stored width is 13.4101845877511 characters. The host's ColumnWidth
# setter quantizes the stored width to (round(chars*6)+5)/6, so the closest
# reachable stored value is 13.333333333333334 (chars input 12.5).
$ws1.Columns.Item(5).ColumnWidth = 12.5   # column E
$ws1.Columns.Item(6).ColumnWidth = 12.5   # column F
$ws2.Columns.Item(3).ColumnWidth = 12.5   # column C
$ws3.Columns.Item(3).ColumnWidth = 12.5   # column C
